$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20
$ws.Range("G20").Value = 2.95
$ws.Range("H20").Value = 2.77
$ws.Range("I20").Value = 2.85
$ws.Range("J20").Value = 3.7
$ws.Range("L20").Value = 3.6
$ws.Range("P20").Value = 2.4
$ws.Range("Q20").Value = 2.72
$ws.Range("U20").Value = 2.07
$ws.Range("V20").Value = 1.7
$ws.Range("W20").Value = 6.7
$ws.Range("X20").Value = 14
$ws.Range("Y20").Value = 11.75
$ws.Range("Z20").Value = 40
$ws.Range("AA20").Value = 35
$ws.Range("AD20").Value = 5.7
$ws.Range("AE20").Value = 18.5
$ws.Range("AH20").Value = 6.6
$ws.Range("AI20").Value = 13.5
$ws.Range("AJ20").Value = 11.5
$ws.Range("AK20").Value = 40
$ws.Range("AM20").Value = 55
$ws.Range("AN20").Value = 4.6
$ws.Range("AO20").Value = 17.5
$ws.Range("AQ20").Value = 90
$ws.Range("AS20").Value = 500
$ws.Range("AU20").Value = 7.6
$ws.Range("AV20").Value = 90
$ws.Range("AY20").Value = 4.45
$ws.Range("AZ20").Value = 16.5
$ws.Range("BA20").Value = 28
$ws.Range("BB20").Value = 90

# Row 37
$ws.Range("G37").Value = 2.2
$ws.Range("I37").Value = 3.2
$ws.Range("S37").Value = 1.36
$ws.Range("W37").Value = 9
$ws.Range("AL37").Value = 23
$ws.Range("AM37").Value = 29

# Row 38
$ws.Range("G38").Value = 2
$ws.Range("I38").Value = 3.9
$ws.Range("J38").Value = 2.63
$ws.Range("L38").Value = 4.5
$ws.Range("M38").Value = 1.07
$ws.Range("N38").Value = 9
$ws.Range("O38").Value = 1.36
$ws.Range("P38").Value = 3.2
$ws.Range("S38").Value = 1.44
$ws.Range("T38").Value = 2.63
$ws.Range("X38").Value = 8.5
$ws.Range("Z38").Value = 17
$ws.Range("AA38").Value = 17
$ws.Range("AB38").Value = 29
$ws.Range("AD38").Value = 6
$ws.Range("AE38").Value = 15
$ws.Range("AH38").Value = 10
$ws.Range("AI38").Value = 19
$ws.Range("AJ38").Value = 15
$ws.Range("AO38").Value = 11
$ws.Range("AT38").Value = 2.63
$ws.Range("AZ38").Value = 23

# Row 41
$ws.Range("H41").Value = 2.9
$ws.Range("M41").Value = 1.11
$ws.Range("N41").Value = 6.5
$ws.Range("O41").Value = 1.5
$ws.Range("P41").Value = 2.5
$ws.Range("Q41").Value = 2.6
$ws.Range("R41").Value = 1.48
$ws.Range("S41").Value = 1.57
$ws.Range("T41").Value = 2.25
$ws.Range("U41").Value = 2.2
$ws.Range("V41").Value = 1.62
$ws.Range("W41").Value = 6
$ws.Range("AA41").Value = 21
$ws.Range("AC41").Value = 6
$ws.Range("AF41").Value = 67
$ws.Range("AH41").Value = 8.5
$ws.Range("AQ41").Value = 41
$ws.Range("AT41").Value = 2.25
$ws.Range("BA41").Value = 34
$ws.Range("BC41").Value = 126

# Row 42
$ws.Range("Q42").Value = 2.4
$ws.Range("R42").Value = 1.53

# Row 74
$ws.Range("G74").Value = 1.6
$ws.Range("H74").Value = 3.5
$ws.Range("I74").Value = 5.6
$ws.Range("J74").Value = 2.12
$ws.Range("M74").Value = 1.05
$ws.Range("N74").Value = 8.85
$ws.Range("P74").Value = 2.72
$ws.Range("S74").Value = 1.4
$ws.Range("T74").Value = 2.5
$ws.Range("U74").Value = 2
$ws.Range("V74").Value = 1.65
$ws.Range("W74").Value = 5.6
$ws.Range("X74").Value = 6.7
$ws.Range("Z74").Value = 11.5
$ws.Range("AB74").Value = 32
$ws.Range("AC74").Value = 8.25
$ws.Range("AD74").Value = 6.9
$ws.Range("AE74").Value = 19
$ws.Range("AH74").Value = 13
$ws.Range("AI74").Value = 32
$ws.Range("AM74").Value = 75
$ws.Range("AN74").Value = 3.3
$ws.Range("AP74").Value = 18.5
$ws.Range("AR74").Value = 60
$ws.Range("AS74").Value = 250
$ws.Range("AT74").Value = 2.47
$ws.Range("AU74").Value = 7.9

# Row 96
$ws.Range("G96").Value = 2.25
$ws.Range("H96").Value = 3.25
$ws.Range("I96").Value = 3.25
$ws.Range("J96").Value = 2.88
$ws.Range("L96").Value = 3.6
$ws.Range("O96").Value = 1.22
$ws.Range("P96").Value = 4
$ws.Range("Q96").Value = 1.8
$ws.Range("R96").Value = 2.05
$ws.Range("U96").Value = 1.62
$ws.Range("V96").Value = 2.2
$ws.Range("W96").Value = 9.5
$ws.Range("X96").Value = 12
$ws.Range("AC96").Value = 12
$ws.Range("AI96").Value = 17
$ws.Range("AL96").Value = 23
$ws.Range("AO96").Value = 12
$ws.Range("AY96").Value = 5

# Row 107
$ws.Range("G107").Value = 3.6
$ws.Range("H107").Value = 2.77
$ws.Range("I107").Value = 2.12
$ws.Range("J107").Value = 4.05
$ws.Range("K107").Value = 1.95
$ws.Range("L107").Value = 2.67
$ws.Range("M107").Value = 1.06
$ws.Range("N107").Value = 8
$ws.Range("O107").Value = 1.36
$ws.Range("P107").Value = 3
$ws.Range("Q107").Value = 2.02
$ws.Range("R107").Value = 1.62
$ws.Range("S107").Value = 1.39
$ws.Range("T107").Value = 2.42
$ws.Range("U107").Value = 1.93
$ws.Range("V107").Value = 1.82
$ws.Range("W107").Value = 8
$ws.Range("X107").Value = 16
$ws.Range("Y107").Value = 10
$ws.Range("Z107").Value = 45
$ws.Range("AA107").Value = 28
$ws.Range("AB107").Value = 30
$ws.Range("AC107").Value = 7.4
$ws.Range("AD107").Value = 4.85
$ws.Range("AE107").Value = 10.75
$ws.Range("AF107").Value = 45
$ws.Range("AH107").Value = 5.9
$ws.Range("AI107").Value = 8.75
$ws.Range("AJ107").Value = 7.1
$ws.Range("AK107").Value = 17.5
$ws.Range("AL107").Value = 14.5
$ws.Range("AM107").Value = 21
$ws.Range("AN107").Value = 5.6
$ws.Range("AO107").Value = 21
$ws.Range("AP107").Value = 25
$ws.Range("AQ107").Value = 110
$ws.Range("AR107").Value = 150
$ws.Range("AS107").Value = 300
$ws.Range("AT107").Value = 2.42
$ws.Range("AU107").Value = 6.3
$ws.Range("AV107").Value = 50
$ws.Range("AY107").Value = 4.05
$ws.Range("AZ107").Value = 11
$ws.Range("BA107").Value = 17.5
$ws.Range("BB107").Value = 45
$ws.Range("BC107").Value = 65

# Row 136
$ws.Range("G136").Value = 7.5
$ws.Range("H136").Value = 4.75
$ws.Range("I136").Value = 1.4
$ws.Range("J136").Value = 8.5
$ws.Range("L136").Value = 2
$ws.Range("U136").Value = 2.63
$ws.Range("V136").Value = 1.44
$ws.Range("Y136").Value = 23
$ws.Range("Z136").Value = 101
$ws.Range("AA136").Value = 67
$ws.Range("AD136").Value = 9.5
$ws.Range("AF136").Value = 126
$ws.Range("AK136").Value = 8.5
$ws.Range("AN136").Value = 9.5
$ws.Range("AQ136").Value = 251
$ws.Range("AR136").Value = 301
$ws.Range("AY136").Value = 3.1
$ws.Range("AZ136").Value = 7
$ws.Range("BB136").Value = 21
